# Generate Report for handback
# Adds two new handback rows (73720230-e572-4bd2-a488-24f546212a77 and
# ed69dab0-1987-48f0-beaf-0f2d634f5c5b) to the Overview, zh-cn and de-de
# sheets of the handback-status workbook, mirroring the pattern used by the
# existing rows.

$wb = $excel.ActiveWorkbook

$uuid1 = "73720230-e572-4bd2-a488-24f546212a77"
$hash1 = "f15c5b67e58affee3dd656f378367ac8ff832b82"
$uuid2 = "ed69dab0-1987-48f0-beaf-0f2d634f5c5b"
$hash2 = "4022d02fe9da449567e8047da85d57405d070d25"

$status = "Handed back: in sync with en-US"
$reason = "Include"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Add($wsOverview.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/" + $uuid1 + ".md", "", "", $uuid1 + ".md")
$wsOverview.Range("B6").Value = $status
$wsOverview.Range("C6").Value = $status

$wsOverview.Hyperlinks.Add($wsOverview.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/" + $uuid2 + ".md", "", "", $uuid2 + ".md")
$wsOverview.Range("B7").Value = $status
$wsOverview.Range("C7").Value = $status

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$xlf1zh = $uuid1 + "." + $hash1 + ".zh-cn.xlf"
$xlf2zh = $uuid2 + "." + $hash2 + ".zh-cn.xlf"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/" + $uuid1 + ".md", "", "", $uuid1 + ".md")
$wsZhCn.Range("B6").Value = $status
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/" + $xlf1zh, "", "", $xlf1zh)
$wsZhCn.Range("D6").Value = "2016-02-15 08:54:25"
$wsZhCn.Range("D6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("E6"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/master/e2e/" + $uuid1 + ".md", "", "", $uuid1 + ".md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F6"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/master/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/" + $xlf1zh, "", "", $xlf1zh)
$wsZhCn.Range("G6").Value = "2016-02-15 08:55:30"
$wsZhCn.Range("G6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("H6").Value = $reason

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/" + $uuid2 + ".md", "", "", $uuid2 + ".md")
$wsZhCn.Range("B7").Value = $status
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/" + $xlf2zh, "", "", $xlf2zh)
$wsZhCn.Range("D7").Value = "2016-02-15 08:54:25"
$wsZhCn.Range("D7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("E7"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/master/e2e/" + $uuid2 + ".md", "", "", $uuid2 + ".md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F7"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/master/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/" + $xlf2zh, "", "", $xlf2zh)
$wsZhCn.Range("G7").Value = "2016-02-15 08:55:30"
$wsZhCn.Range("G7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("H7").Value = $reason

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$xlf1de = $uuid1 + "." + $hash1 + ".de-de.xlf"
$xlf2de = $uuid2 + "." + $hash2 + ".de-de.xlf"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/" + $uuid1 + ".md", "", "", $uuid1 + ".md")
$wsDeDe.Range("B6").Value = $status
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/" + $xlf1de, "", "", $xlf1de)
$wsDeDe.Range("D6").Value = "2016-02-15 08:54:44"
$wsDeDe.Range("D6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("E6"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/master/e2e/" + $uuid1 + ".md", "", "", $uuid1 + ".md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F6"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/master/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/" + $xlf1de, "", "", $xlf1de)
$wsDeDe.Range("G6").Value = "2016-02-15 08:55:59"
$wsDeDe.Range("G6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("H6").Value = $reason

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/" + $uuid2 + ".md", "", "", $uuid2 + ".md")
$wsDeDe.Range("B7").Value = $status
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/" + $xlf2de, "", "", $xlf2de)
$wsDeDe.Range("D7").Value = "2016-02-15 08:54:44"
$wsDeDe.Range("D7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("E7"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/master/e2e/" + $uuid2 + ".md", "", "", $uuid2 + ".md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F7"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/master/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/" + $xlf2de, "", "", $xlf2de)
$wsDeDe.Range("G7").Value = "2016-02-15 08:55:59"
$wsDeDe.Range("G7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("H7").Value = $reason

Write-Host "Handback rows added for" $uuid1 "and" $uuid2
